$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2384615384615385
$ws.Range("C2").Value = 0.4692307692307692
$ws.Range("J2").Value = 0.01153846153846154
$ws.Range("P2").Value = 0.1846153846153846
$ws.Range("S2").Value = 0.09615384615384616
$ws.Range("B3").Value = 0.008130081300813009
$ws.Range("C3").Value = 0.01626016260162602
$ws.Range("P3").Value = 0.7235772357723578
$ws.Range("S3").Value = 0.2520325203252032
$ws.Range("J4").Value = 0.06
$ws.Range("P4").Value = 0.66
$ws.Range("S4").Value = 0.28
$ws.Range("B6").Value = 0.07981220657276995
$ws.Range("D6").Value = 0.01408450704225352
$ws.Range("F6").Value = 0.03286384976525822
$ws.Range("J6").Value = 0.3098591549295774
$ws.Range("O6").Value = 0.01408450704225352
$ws.Range("Q6").Value = 0.1690140845070423
$ws.Range("R6").Value = 0.04694835680751173
$ws.Range("S6").Value = 0.3333333333333333
$ws.Range("B7").Value = 0.1006289308176101
$ws.Range("D7").Value = 0.01257861635220126
$ws.Range("E7").Value = 0.006289308176100629
$ws.Range("F7").Value = 0.07547169811320754
$ws.Range("J7").Value = 0.1069182389937107
$ws.Range("O7").Value = 0.01886792452830189
$ws.Range("Q7").Value = 0.1886792452830189
$ws.Range("R7").Value = 0.1257861635220126
$ws.Range("S7").Value = 0.3647798742138365
$ws.Range("B8").Value = 0.06094808126410835
$ws.Range("D8").Value = 0.01580135440180587
$ws.Range("E8").Value = 0.006772009029345372
$ws.Range("F8").Value = 0.06094808126410835
$ws.Range("J8").Value = 0.1376975169300226
$ws.Range("O8").Value = 0.02708803611738149
$ws.Range("Q8").Value = 0.1918735891647856
$ws.Range("R8").Value = 0.1151241534988713
$ws.Range("S8").Value = 0.3837471783295711
$ws.Range("B9").Value = 0.08421052631578947
$ws.Range("D9").Value = 0.02105263157894737
$ws.Range("F9").Value = 0.07894736842105263
$ws.Range("J9").Value = 0.1
$ws.Range("O9").Value = 0.01052631578947368
$ws.Range("Q9").Value = 0.2
$ws.Range("R9").Value = 0.131578947368421
$ws.Range("S9").Value = 0.3736842105263158
$ws.Range("B10").Value = 0.09456264775413711
$ws.Range("D10").Value = 0.02758077226162333
$ws.Range("F10").Value = 0.05752561071710008
$ws.Range("J10").Value = 0.1268715524034673
$ws.Range("O10").Value = 0.01339637509850276
$ws.Range("Q10").Value = 0.22301024428684
$ws.Range("R10").Value = 0.08274231678486997
$ws.Range("S10").Value = 0.3743104806934594
$ws.Range("G11").Value = 0.1488549618320611
$ws.Range("J11").Value = 0.1068702290076336
$ws.Range("K11").Value = 0.2175572519083969
$ws.Range("L11").Value = 0.5152671755725191
$ws.Range("S11").Value = 0.01145038167938931
$ws.Range("G12").Value = 0.696551724137931
$ws.Range("J12").Value = 0.2137931034482759
$ws.Range("K12").Value = 0.01379310344827586
$ws.Range("L12").Value = 0.04827586206896552
$ws.Range("S12").Value = 0.02758620689655172
$ws.Range("G13").Value = 0.5434782608695652
$ws.Range("J13").Value = 0.4565217391304348
$ws.Range("F15").Value = 0.02643171806167401
$ws.Range("H15").Value = 0.1629955947136564
$ws.Range("I15").Value = 0.07488986784140969
$ws.Range("J15").Value = 0.4096916299559472
$ws.Range("K15").Value = 0.06167400881057269
$ws.Range("M15").Value = 0.013215859030837
$ws.Range("N15").Value = 0.004405286343612335
$ws.Range("O15").Value = 0.04405286343612335
$ws.Range("S15").Value = 0.2026431718061674
$ws.Range("F16").Value = 0.01785714285714286
$ws.Range("I16").Value = 0.130952380952381
$ws.Range("J16").Value = 0.375
$ws.Range("K16").Value = 0.07738095238095238
$ws.Range("M16").Value = 0.01785714285714286
$ws.Range("O16").Value = 0.06547619047619048
$ws.Range("S16").Value = 0.125
$ws.Range("F17").Value = 0.02553191489361702
$ws.Range("H17").Value = 0.2021276595744681
$ws.Range("I17").Value = 0.07659574468085106
$ws.Range("J17").Value = 0.4404255319148936
$ws.Range("K17").Value = 0.08936170212765958
$ws.Range("M17").Value = 0.01702127659574468
$ws.Range("N17").Value = 0.002127659574468085
$ws.Range("O17").Value = 0.05957446808510639
$ws.Range("S17").Value = 0.08723404255319149
$ws.Range("F18").Value = 0.02380952380952381
$ws.Range("H18").Value = 0.1809523809523809
$ws.Range("I18").Value = 0.0761904761904762
$ws.Range("J18").Value = 0.3952380952380952
$ws.Range("K18").Value = 0.119047619047619
$ws.Range("M18").Value = 0.01428571428571429
$ws.Range("N18").Value = 0.004761904761904762
$ws.Range("O18").Value = 0.0761904761904762
$ws.Range("S18").Value = 0.1095238095238095
$ws.Range("F19").Value = 0.02464788732394366
$ws.Range("H19").Value = 0.2121478873239437
$ws.Range("I19").Value = 0.08714788732394366
$ws.Range("J19").Value = 0.3785211267605634
$ws.Range("K19").Value = 0.09595070422535211
$ws.Range("M19").Value = 0.02464788732394366
$ws.Range("N19").Value = 0.00176056338028169
$ws.Range("O19").Value = 0.08450704225352113
$ws.Range("S19").Value = 0.09066901408450705
